$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51: date moved to 2021-11-09, quality changed Extra -> Primera, new prices ---
$ws.Cells.Item(51, 4).Value2 = 44509
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value2 = 50
$ws.Cells.Item(51, 11).Value2 = 1200
$ws.Cells.Item(51, 12).Value2 = 1200
$ws.Cells.Item(51, 13).Value2 = 1200
$ws.Cells.Item(51, 16).Value2 = 1200

# --- Row 52: quality changed Primera -> Extra, new prices (date stays 2021-10-22) ---
$ws.Cells.Item(52, 9).Value = "Extra"
$ws.Cells.Item(52, 10).Value2 = 250
$ws.Cells.Item(52, 11).Value2 = 1500
$ws.Cells.Item(52, 12).Value2 = 1500
$ws.Cells.Item(52, 13).Value2 = 1500
$ws.Cells.Item(52, 16).Value2 = 1500

# --- Row 53: quality changed Segunda -> Primera, new prices (date stays 2021-10-22) ---
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value2 = 400
$ws.Cells.Item(53, 11).Value2 = 1300
$ws.Cells.Item(53, 12).Value2 = 1300
$ws.Cells.Item(53, 13).Value2 = 1300
$ws.Cells.Item(53, 16).Value2 = 1300

# --- Row 54 (new row): carries over the old "Segunda" record that used to live on row 53 ---
$ws.Cells.Item(54, 1).Value2 = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value2 = 44491
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value2 = 9
$ws.Cells.Item(54, 6).Value2 = 300000000
$ws.Cells.Item(54, 7).Value = "Espárragos"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Segunda"
$ws.Cells.Item(54, 10).Value2 = 300
$ws.Cells.Item(54, 11).Value2 = 1000
$ws.Cells.Item(54, 12).Value2 = 1000
$ws.Cells.Item(54, 13).Value2 = 1000
$ws.Cells.Item(54, 14).Value = "$/kilo"
$ws.Cells.Item(54, 15).Value = "Región del Maule"
$ws.Cells.Item(54, 16).Value2 = 1000
$ws.Cells.Item(54, 17).Value2 = 1
$ws.Cells.Item(54, 18).Value = "Hortaliza"
